# Add an "Up"/"Down" sentiment category column and a new trade row (Long/Hold trades)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: add new columns X (numeric) and Y (shared string "Up") ---
$ws.Range("X3").Value = 0.069999999999993179
$ws.Range("Y3").Value = "Up"

# --- Row 4: new row, copy number formats from row 3 first so styles match (date + percent) ---
$ws.Range("A3:W3").Copy()
$ws.Range("A4:W4").PasteSpecial(-4122)

$ws.Range("A4").Value = 42633.884317129632
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "Neutral"
$ws.Range("D4").Value = 18
$ws.Range("E4").Value = 30090
$ws.Range("F4").Value = 3119
$ws.Range("G4").Value = 57
$ws.Range("H4").Value = 36
$ws.Range("I4").Value = 85
$ws.Range("J4").Value = 15
$ws.Range("K4").Value = 14002
$ws.Range("L4").Value = 386
$ws.Range("M4").Value = 242
$ws.Range("N4").Value = 34
$ws.Range("O4").Value = 6
$ws.Range("P4").Value = "Bag"
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = -31.95
$ws.Range("S4").Value = -0.082
$ws.Range("T4").Value = -0.28000000000000003
$ws.Range("U4").Value = 6.77
$ws.Range("V4").Value = 1.88
$ws.Range("W4").Value = 0

Write-Output "edit complete"
